# Updated GBDS File | August 2025
# Applies the August 2025 purchase entries to the "PE, DECEMBER" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PE, DECEMBER")

# Row 8 - new purchase entry
$ws.Range("C8").Value = 45885
$ws.Range("G8").Value = 517818948
$ws.Range("I8").Formula = "=1366324-57713.04"

# Row 9 - new purchase entry
$ws.Range("C9").Value = 45885
$ws.Range("G9").Value = 517819044
$ws.Range("I9").Formula = "=1353132-56255.04"

# Row 10 - new purchase entry
$ws.Range("C10").Value = 45887
$ws.Range("G10").Value = 517821628
$ws.Range("I10").Formula = "=1353132-56255.04"

# Update the active selection to reflect where the user left off editing
$ws.Activate()
$ws.Range("I10").Select()
